# Insert a new weekly record at row 189 of the "Poroto verde" sheet.
# This pushes the previous rows 189-193 down to 190-194 (unchanged),
# and populates the new row 189 with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 189, shifting rows 189:193 down to 190:194.
$ws.Rows.Item(189).Insert()

# Populate the newly inserted row 189 with the new record.
$ws.Range("A189").Value = 8
$ws.Range("B189").Value = "Terminal La Palmera de La Serena"
$ws.Range("C189").Value = "Coquimbo"
$ws.Range("D189").Value = 44628
$ws.Range("E189").Value = 4
$ws.Range("F189").Value = 100112031
$ws.Range("G189").Value = "Poroto verde"
$ws.Range("H189").Value = "Magnum"
$ws.Range("I189").Value = "Primera"
$ws.Range("J189").Value = 500
$ws.Range("K189").Value = 24000
$ws.Range("L189").Value = 25000
$ws.Range("M189").Value = 24500
$ws.Range("N189").Value = "`$/malla 25 kilos"
$ws.Range("O189").Value = "Provincia de Limarí"
$ws.Range("P189").Value = 980
$ws.Range("Q189").Value = 25
$ws.Range("R189").Value = "Hortaliza"
